$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 8888
$ws.Range("E2").Value = -396
$ws.Range("F2").Value = -396
$ws.Range("G2").Value = -516
$ws.Range("H2").Value = -422
$ws.Range("I2").Value = -422
$ws.Range("K2").Value = 15403
$ws.Range("L2").Value = 8396
$ws.Range("M2").Value = 7007
$ws.Range("N2").Value = 7007
$ws.Range("P2").Value = 695
$ws.Range("Q2").Value = -855
$ws.Range("R2").Value = 964
$ws.Range("S2").Value = -77
$ws.Range("T2").Value = 18
$ws.Range("U2").Value = -873
$ws.Range("V2").Value = 3043
$ws.Range("W2").Value = -4.46
$ws.Range("X2").Value = -4.75
$ws.Range("Y2").Value = -5.67
$ws.Range("Z2").Value = -2.64
$ws.Range("AA2").Value = 119.82
$ws.Range("AB2").Value = 955.9400000000001
$ws.Range("AC2").Value = -608
$ws.Range("AE2").Value = 10082
$ws.Range("AF2").Value = 0.6
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 69500000
$ws.Range("J2").ClearContents()
$ws.Range("O2").ClearContents()
$ws.Range("AD2").ClearContents()

# Row 3
$ws.Range("D3").Value = 6936
$ws.Range("E3").Value = -638
$ws.Range("F3").Value = -638
$ws.Range("G3").Value = -1657
$ws.Range("H3").Value = -1254
$ws.Range("I3").Value = -1254
$ws.Range("K3").Value = 14060
$ws.Range("L3").Value = 8314
$ws.Range("M3").Value = 5747
$ws.Range("N3").Value = 5747
$ws.Range("P3").Value = 695
$ws.Range("Q3").Value = -441
$ws.Range("R3").Value = 149
$ws.Range("S3").Value = 296
$ws.Range("T3").Value = 50
$ws.Range("U3").Value = -492
$ws.Range("V3").Value = 3346
$ws.Range("W3").Value = -9.19
$ws.Range("X3").Value = -18.09
$ws.Range("Y3").Value = -19.67
$ws.Range("Z3").Value = -8.52
$ws.Range("AA3").Value = 144.67
$ws.Range("AB3").Value = 779.76
$ws.Range("AC3").Value = -1805
$ws.Range("AE3").Value = 8269
$ws.Range("AF3").Value = 0.37
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 69500000
$ws.Range("J3").ClearContents()
$ws.Range("O3").ClearContents()
$ws.Range("AD3").ClearContents()

# Row 4
$ws.Range("D4").Value = 8029
$ws.Range("E4").Value = 42
$ws.Range("F4").Value = 42
$ws.Range("G4").Value = -1878
$ws.Range("H4").Value = -1812
$ws.Range("I4").Value = -1812
$ws.Range("K4").Value = 13564
$ws.Range("L4").Value = 8222
$ws.Range("M4").Value = 5343
$ws.Range("N4").Value = 5343
$ws.Range("P4").Value = 695
$ws.Range("Q4").Value = 45
$ws.Range("R4").Value = 414
$ws.Range("S4").Value = -187
$ws.Range("T4").Value = 168
$ws.Range("U4").Value = -123
$ws.Range("V4").Value = 3162
$ws.Range("W4").Value = 0.53
$ws.Range("X4").Value = -22.57
$ws.Range("Y4").Value = -32.68
$ws.Range("Z4").Value = -13.12
$ws.Range("AA4").Value = 153.89
$ws.Range("AB4").Value = 525.2
$ws.Range("AC4").Value = -2608
$ws.Range("AE4").Value = 7687
$ws.Range("AF4").Value = 0.4
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 69500000
$ws.Range("J4").ClearContents()
$ws.Range("O4").ClearContents()
$ws.Range("AD4").ClearContents()

# Row 5
$ws.Range("D5").Value = 7689
$ws.Range("E5").Value = 135
$ws.Range("F5").Value = 135
$ws.Range("G5").Value = -49
$ws.Range("H5").Value = -103
$ws.Range("I5").Value = -103
$ws.Range("K5").Value = 12313
$ws.Range("L5").Value = 6882
$ws.Range("M5").Value = 5432
$ws.Range("N5").Value = 5432
$ws.Range("P5").Value = 695
$ws.Range("Q5").Value = -694
$ws.Range("R5").Value = 127
$ws.Range("S5").Value = 566
$ws.Range("T5").Value = 29
$ws.Range("U5").Value = -723
$ws.Range("V5").Value = 3731
$ws.Range("W5").Value = 1.75
$ws.Range("X5").Value = -1.35
$ws.Range("Y5").Value = -1.92
$ws.Range("Z5").Value = -0.8
$ws.Range("AA5").Value = 126.7
$ws.Range("AB5").Value = 514.09
$ws.Range("AC5").Value = -149
$ws.Range("AD5").Value = -23.41
$ws.Range("AE5").Value = 7815
$ws.Range("AF5").Value = 0.45
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 69500000
$ws.Range("J5").ClearContents()
$ws.Range("O5").ClearContents()

# Row 6
$ws.Range("D6").Value = 5113
$ws.Range("E6").Value = -353
$ws.Range("F6").Value = -353
$ws.Range("G6").Value = -296
$ws.Range("H6").Value = -187
$ws.Range("I6").Value = -187
$ws.Range("K6").Value = 7434
$ws.Range("L6").Value = 5124
$ws.Range("M6").Value = 2310
$ws.Range("N6").Value = 2310
$ws.Range("P6").Value = 329
$ws.Range("Q6").Value = 184
$ws.Range("R6").Value = -148
$ws.Range("S6").Value = -161
$ws.Range("T6").Value = 25
$ws.Range("U6").Value = 160
$ws.Range("V6").Value = 2236
$ws.Range("W6").Value = -6.91
$ws.Range("X6").Value = -3.65
$ws.Range("Y6").Value = -4.83
$ws.Range("Z6").Value = -1.89
$ws.Range("AA6").Value = 221.78
$ws.Range("AB6").Value = 548.42
$ws.Range("AC6").Value = -385
$ws.Range("AD6").Value = -13.7
$ws.Range("AE6").Value = 7014
$ws.Range("AF6").Value = 0.75
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 32947142
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()

# Row 7
$ws.Range("D7").Value = 6387
$ws.Range("E7").Value = -150
$ws.Range("G7").Value = -278
$ws.Range("H7").Value = -236
$ws.Range("I7").Value = -236
$ws.Range("K7").Value = 8547
$ws.Range("L7").Value = 6049
$ws.Range("M7").Value = 2088
$ws.Range("N7").Value = 2088
$ws.Range("P7").Value = 330
$ws.Range("Q7").Value = -170
$ws.Range("R7").Value = -196
$ws.Range("S7").Value = 605
$ws.Range("T7").Value = 144
$ws.Range("U7").Value = 705
$ws.Range("W7").Value = -2.35
$ws.Range("X7").Value = -3.69
$ws.Range("Y7").Value = -10.72
$ws.Range("Z7").Value = -2.95
$ws.Range("AA7").Value = 289.65
$ws.Range("AC7").Value = -716
$ws.Range("AD7").Value = -5.67
$ws.Range("AE7").Value = 6340
$ws.Range("AF7").Value = 0.64
$ws.Range("AG7").Value = 0
$ws.Range("AH7").Value = 0
$ws.Range("AI7").ClearContents()

# Row 8
$ws.Range("D8").Value = 8117
$ws.Range("E8").Value = 230
$ws.Range("G8").Value = 102
$ws.Range("H8").Value = 81
$ws.Range("I8").Value = 81
$ws.Range("K8").Value = 9546
$ws.Range("L8").Value = 6582
$ws.Range("M8").Value = 2173
$ws.Range("N8").Value = 2173
$ws.Range("P8").Value = 330
$ws.Range("Q8").Value = 183
$ws.Range("R8").Value = -190
$ws.Range("S8").Value = 270
$ws.Range("T8").Value = 132
$ws.Range("U8").Value = 650
$ws.Range("W8").Value = 2.83
$ws.Range("X8").Value = 1
$ws.Range("Y8").Value = 3.81
$ws.Range("Z8").Value = 0.9
$ws.Range("AA8").Value = 302.85
$ws.Range("AC8").Value = 247
$ws.Range("AD8").Value = 16.46
$ws.Range("AE8").Value = 6598
$ws.Range("AF8").Value = 0.62
$ws.Range("AG8").Value = 0
$ws.Range("AH8").Value = 0
$ws.Range("AI8").ClearContents()

# Row 9
$ws.Range("D9").Value = 8684
$ws.Range("E9").Value = 322
$ws.Range("G9").Value = 186
$ws.Range("H9").Value = 147
$ws.Range("I9").Value = 147
$ws.Range("K9").Value = 10049
$ws.Range("L9").Value = 6708
$ws.Range("M9").Value = 2325
$ws.Range("N9").Value = 2325
$ws.Range("P9").Value = 330
$ws.Range("Q9").Value = 460
$ws.Range("R9").Value = -153
$ws.Range("S9").Value = -16
$ws.Range("T9").Value = 148
$ws.Range("U9").Value = 514
$ws.Range("W9").Value = 3.7
$ws.Range("X9").Value = 1.69
$ws.Range("Y9").Value = 6.53
$ws.Range("Z9").Value = 1.5
$ws.Range("AA9").Value = 288.49
$ws.Range("AC9").Value = 445
$ws.Range("AD9").Value = 9.119999999999999
$ws.Range("AE9").Value = 7058
$ws.Range("AF9").Value = 0.58
$ws.Range("AG9").Value = 0
$ws.Range("AH9").Value = 0
$ws.Range("AI9").ClearContents()
